# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole "Price" column to Text format so that values such as
# "1.003" or "12.20" are written back as literal strings (matching the
# original inlineStr cells) instead of being auto-coerced into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "28.106.92"
$ws.Range("E2").Value = "  -1.43%  "

# Row 3
$ws.Range("D3").Value = "1.793.21"
$ws.Range("E3").Value = "  -1.59%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "314.09"
$ws.Range("E5").Value = "  -0.38%  "

# Row 6
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.17%  "

# Row 7
$ws.Range("D7").Value = "0.5222"
$ws.Range("E7").Value = "  +2.17%  "

# Row 8
$ws.Range("D8").Value = "0.3813"
$ws.Range("E8").Value = "  -3.46%  "

# Row 9
$ws.Range("D9").Value = "0.07950"
$ws.Range("E9").Value = "  -2.63%  "

# Row 10
$ws.Range("D10").Value = "41.42"
$ws.Range("E10").Value = "  -0.53%  "

# Row 11
$ws.Range("D11").Value = "1.092"
$ws.Range("E11").Value = "  -1.68%  "

# Row 12
$ws.Range("D12").Value = "6.272"
$ws.Range("E12").Value = "  -0.70%  "

# Row 13
$ws.Range("E13").Value = "  +0.14%  "

# Row 14
$ws.Range("D14").Value = "20.47"
$ws.Range("E14").Value = "  -3.01%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.795.50"
$ws.Range("E15").Value = "  -1.55%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.287"
$ws.Range("E16").Value = "  -3.06%  "

# Row 17
$ws.Range("D17").Value = "91.82"
$ws.Range("E17").Value = "  -0.81%  "

# Row 18
$ws.Range("D18").Value = "0.00001083"
$ws.Range("E18").Value = "  -3.97%  "

# Row 19
$ws.Range("D19").Value = "0.06570"
$ws.Range("E19").Value = "  -1.28%  "

# Row 20
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("D21").Value = "17.29"
$ws.Range("E21").Value = "  -2.73%  "

# Row 22
$ws.Range("D22").Value = "5.940"
$ws.Range("E22").Value = "  -2.46%  "

# Row 23
$ws.Range("D23").Value = "28.143.75"
$ws.Range("E23").Value = "  -1.44%  "

# Row 24
$ws.Range("D24").Value = "11.12"
$ws.Range("E24").Value = "  -2.41%  "

# Row 26
$ws.Range("D26").Value = "160.60"
$ws.Range("E26").Value = "  +2.62%  "

# Row 27
$ws.Range("E27").Value = "  -4.77%  "

# Row 28
$ws.Range("D28").Value = "1.998.13"
$ws.Range("E28").Value = "  -1.72%  "

# Row 29
$ws.Range("D29").Value = "2.328"
$ws.Range("E29").Value = "  -3.11%  "

# Row 30
$ws.Range("D30").Value = "122.69"
$ws.Range("E30").Value = "  -2.58%  "

# Row 31
$ws.Range("D31").Value = "0.1078"
$ws.Range("E31").Value = "  -1.25%  "

# Row 32
$ws.Range("D32").Value = "1.051"
$ws.Range("E32").Value = "  -5.60%  "

# Row 33
$ws.Range("D33").Value = "3.673"
$ws.Range("E33").Value = "  +0.31%  "

# Row 34
$ws.Range("D34").Value = "5.526"
$ws.Range("E34").Value = "  -3.97%  "

# Row 35
$ws.Range("D35").Value = "0.07219"
$ws.Range("E35").Value = "  +2.33%  "

# Row 36
$ws.Range("D36").Value = "12.20"
$ws.Range("E36").Value = "  +8.09%  "

# Row 37
$ws.Range("D37").Value = "0.02317"
$ws.Range("E37").Value = "  -1.40%  "

# Row 38
$ws.Range("D38").Value = "8.789"
$ws.Range("E38").Value = "  -0.55%  "

# Row 39
$ws.Range("D39").Value = "0.2139"
$ws.Range("E39").Value = "  -4.03%  "

# Row 40
$ws.Range("D40").Value = "5.059"
$ws.Range("E40").Value = "  -4.25%  "

# Row 41
$ws.Range("D41").Value = "0.6137"
$ws.Range("E41").Value = "  -2.93%  "

# Row 42
$ws.Range("D42").Value = "1.163"
$ws.Range("E42").Value = "  -1.49%  "

# Row 43
$ws.Range("D43").Value = "1.364"
$ws.Range("E43").Value = "  -2.41%  "

# Row 44
$ws.Range("D44").Value = "13.23"
$ws.Range("E44").Value = "  -1.83%  "

# Row 45
$ws.Range("D45").Value = "3.771"
$ws.Range("E45").Value = "  +1.08%  "

# Row 46
$ws.Range("D46").Value = "0.5967"
$ws.Range("E46").Value = "  +0.66%  "

# Row 47
$ws.Range("D47").Value = "127.60"
$ws.Range("E47").Value = "  +1.93%  "

# Row 48
$ws.Range("D48").Value = "1.230"
$ws.Range("E48").Value = "  +3.47%  "

# Row 49
$ws.Range("D49").Value = "1.913"
$ws.Range("E49").Value = "  -3.86%  "

# Row 50
$ws.Range("D50").Value = "0.06739"
$ws.Range("E50").Value = "  -2.40%  "

# Row 51
$ws.Range("D51").Value = "72.86"
$ws.Range("E51").Value = "  -1.55%  "

# Restore the column to its original (default/General) style now that the
# text values have been safely written, so cell formatting matches the
# untouched cells (e.g. row 25 / Toncoin) exactly.
$priceRange.ClearFormats()
